$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 10 values: B10 22->23, C10 24->23
$ws.Range("B10").Value = 23
$ws.Range("C10").Value = 23

# Remove E10 entirely (it held the shared string "dao074", now deleted)
$ws.Range("E10").Clear()

# Update active selection to E11
$ws.Range("E11").Select()
